$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292-354 down to 293-355.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new observation.
$ws.Range("A292").Value = 4
$ws.Range("B292").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C292").Value = "Los Lagos"
$ws.Range("D292").Value = 44641
$ws.Range("E292").Value = 10
$ws.Range("F292").Value = 100114001
$ws.Range("G292").Value = "Papa"
$ws.Range("H292").Value = "Patagonia"
$ws.Range("I292").Value = "1a (cosecha)"
$ws.Range("J292").Value = 300
$ws.Range("K292").Value = 7000
$ws.Range("L292").Value = 8000
$ws.Range("M292").Value = 7500
$ws.Range("N292").Value = "`$/saco 25 kilos"
$ws.Range("O292").Value = "Provincia de Llanquihue"
$ws.Range("P292").Value = 300
$ws.Range("Q292").Value = 25
$ws.Range("R292").Value = "Hortaliza"
